$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 7.155837786072367
$ws.Range("G2").Value = 2781
$ws.Range("H2").Value = 5.491031224551755
$ws.Range("I2").Value = 1.303186504217432
$ws.Range("J2").Value = 191.5684161199625
$ws.Range("K2").Value = 1051.908154552638
$ws.Range("L2").Value = 807.181590009108
$ws.Range("M2").Value = 1.187153488709406
$ws.Range("F3").Value = 7.189756796779523
$ws.Range("G3").Value = 2776
$ws.Range("H3").Value = 5.526996039022875
$ws.Range("I3").Value = 1.300843486410497
$ws.Range("J3").Value = 163.9062792877226
$ws.Range("K3").Value = 905.9093563942199
$ws.Range("L3").Value = 696.4015009168822
$ws.Range("M3").Value = 1.299908028857738
$ws.Range("F4").Value = 7.144365214794698
$ws.Range("G4").Value = 2778
$ws.Range("H4").Value = 5.488148080767418
$ws.Range("I4").Value = 1.301780693533271
$ws.Range("J4").Value = 1.301780693533271
$ws.Range("K4").Value = 7.144365214794698
$ws.Range("L4").Value = 5.488148080767418
$ws.Range("M4").Value = 1.786805740220154
$ws.Range("F5").Value = 7.226133001156112
$ws.Range("G5").Value = 2784
$ws.Range("H5").Value = 5.538997063386187
$ws.Range("I5").Value = 1.304592314901593
$ws.Range("J5").Value = 1.304592314901593
$ws.Range("K5").Value = 7.226133001156112
$ws.Range("L5").Value = 5.538997063386187
$ws.Range("M5").Value = 1.902640819204404
$ws.Range("F6").Value = 7.831953062165842
$ws.Range("G6").Value = 2773
$ws.Range("H6").Value = 6.027186381053698
$ws.Range("I6").Value = 1.299437675726336
$ws.Range("J6").Value = 171.5257731958763
$ws.Range("K6").Value = 1033.817804205891
$ws.Range("L6").Value = 795.5886022990882
$ws.Range("M6").Value = 1.621214283868329
$ws.Range("F7").Value = 8.025883211335147
$ws.Range("G7").Value = 2837
$ws.Range("H7").Value = 6.037093680997252
$ws.Range("I7").Value = 1.329428303655108
$ws.Range("J7").Value = 195.4259606373008
$ws.Range("K7").Value = 1179.804832066267
$ws.Range("L7").Value = 887.4527711065961
$ws.Range("M7").Value = 1.331494024760501
$ws.Range("F8").Value = 8.132159920605083
$ws.Range("G8").Value = 2856
$ws.Range("H8").Value = 6.076340781012342
$ws.Range("I8").Value = 1.338331771321462
$ws.Range("J8").Value = 1.338331771321462
$ws.Range("K8").Value = 8.132159920605083
$ws.Range("L8").Value = 6.076340781012342
$ws.Range("M8").Value = 2.486001287728974
$ws.Range("F9").Value = 7.92476738644709
$ws.Range("G9").Value = 2789
$ws.Range("H9").Value = 6.063626246926529
$ws.Range("I9").Value = 1.306935332708529
$ws.Range("J9").Value = 103.2478912839738
$ws.Range("K9").Value = 626.0566235293202
$ws.Range("L9").Value = 479.0264735071958
$ws.Range("M9").Value = 1.231508851853878
$ws.Range("F10").Value = 9.970799587348667
$ws.Range("G10").Value = 2831
$ws.Range("H10").Value = 7.515961257294968
$ws.Range("I10").Value = 1.326616682286785
$ws.Range("J10").Value = 1.326616682286785
$ws.Range("K10").Value = 9.970799587348667
$ws.Range("L10").Value = 7.515961257294968
$ws.Range("M10").Value = 2.822733363178408
$ws.Range("F11").Value = 9.880966773723198
$ws.Range("G11").Value = 2783
$ws.Range("H11").Value = 7.576709700009093
$ws.Range("I11").Value = 1.304123711340206
$ws.Range("J11").Value = 99.11340206185567
$ws.Range("K11").Value = 750.9534748029631
$ws.Range("L11").Value = 575.829937200691
$ws.Range("M11").Value = 1.566133233635127
$ws.Range("F12").Value = 10.26719565501607
$ws.Range("G12").Value = 3092
$ws.Range("H12").Value = 7.086091697220018
$ws.Range("I12").Value = 1.44892221180881
$ws.Range("J12").Value = 114.464854732896
$ws.Range("K12").Value = 811.1084567462696
$ws.Range("L12").Value = 559.8012440803815
$ws.Range("M12").Value = 1.595522204789497
$ws.Range("F13").Value = 10.29119117449056
$ws.Range("G13").Value = 3095
$ws.Range("H13").Value = 7.095768002055853
$ws.Range("I13").Value = 1.450328022492971
$ws.Range("J13").Value = 110.2249297094658
$ws.Range("K13").Value = 782.1305292612828
$ws.Range("L13").Value = 539.2783681562448
$ws.Range("M13").Value = 1.631153801156754
$ws.Range("F14").Value = 11.07777372186156
$ws.Range("G14").Value = 3330
$ws.Range("H14").Value = 7.099089826562332
$ws.Range("I14").Value = 1.560449859418932
$ws.Range("J14").Value = 1.560449859418932
$ws.Range("K14").Value = 11.07777372186156
$ws.Range("L14").Value = 7.099089826562332
$ws.Range("M14").Value = 2.770551207837575
$ws.Range("F15").Value = 10.03223148110968
$ws.Range("G15").Value = 3471
$ws.Range("H15").Value = 6.167900311347755
$ws.Range("I15").Value = 1.626522961574508
$ws.Range("J15").Value = 128.4953139643861
$ws.Range("K15").Value = 792.5462870076648
$ws.Range("L15").Value = 487.2641245964726
$ws.Range("M15").Value = 1.559008772164445
$ws.Range("F16").Value = 10.20465831876525
$ws.Range("G16").Value = 3540
$ws.Range("H16").Value = 6.151621709673738
$ws.Range("I16").Value = 1.658856607310216
$ws.Range("J16").Value = 1.658856607310216
$ws.Range("K16").Value = 10.20465831876525
$ws.Range("L16").Value = 6.151621709673738
$ws.Range("M16").Value = 1.968478589689816
$ws.Range("F17").Value = 9.923597403723614
$ws.Range("G17").Value = 3472
$ws.Range("H17").Value = 6.099353934201092
$ws.Range("I17").Value = 1.626991565135895
$ws.Range("J17").Value = 123.651358950328
$ws.Range("K17").Value = 754.1934026829946
$ws.Range("L17").Value = 463.550898999283
$ws.Range("M17").Value = 1.572890188490193
$ws.Range("F18").Value = 9.938955107279437
$ws.Range("G18").Value = 3470
$ws.Range("H18").Value = 6.112314178367239
$ws.Range("I18").Value = 1.626054358013121
$ws.Range("J18").Value = 1.626054358013121
$ws.Range("K18").Value = 9.938955107279437
$ws.Range("L18").Value = 6.112314178367239
$ws.Range("M18").Value = 2.485732672330587
$ws.Range("F19").Value = 11.05661031699668
$ws.Range("G19").Value = 3774
$ws.Range("H19").Value = 6.251935987406178
$ws.Range("I19").Value = 1.768509840674789
$ws.Range("J19").Value = 1.768509840674789
$ws.Range("K19").Value = 11.05661031699668
$ws.Range("L19").Value = 6.251935987406178
$ws.Range("M19").Value = 2.911205496465226
$ws.Range("F20").Value = 11.58819121916773
$ws.Range("G20").Value = 3472
$ws.Range("H20").Value = 7.122465455559888
$ws.Range("I20").Value = 1.626991565135895
$ws.Range("J20").Value = 128.5323336457357
$ws.Range("K20").Value = 915.4671063142505
$ws.Range("L20").Value = 562.6747709892311
$ws.Range("M20").Value = 1.800804915458665
$ws.Range("F21").Value = 11.85223903773185
$ws.Range("G21").Value = 3540
$ws.Range("H21").Value = 7.144824323875643
$ws.Range("I21").Value = 1.658856607310216
$ws.Range("J21").Value = 1.658856607310216
$ws.Range("K21").Value = 11.85223903773185
$ws.Range("L21").Value = 7.144824323875643
$ws.Range("M21").Value = 2.286296910378474
$ws.Range("F22").Value = 11.64372325453217
$ws.Range("G22").Value = 3467
$ws.Range("H22").Value = 7.166918207433414
$ws.Range("I22").Value = 1.62464854732896
$ws.Range("J22").Value = 1.62464854732896
$ws.Range("K22").Value = 11.64372325453217
$ws.Range("L22").Value = 7.166918207433414
$ws.Range("M22").Value = 3.296338053358058
$ws.Range("F23").Value = 11.53149073933663
$ws.Range("G23").Value = 3464
$ws.Range("H23").Value = 7.103984191034749
$ws.Range("I23").Value = 1.623242736644799
$ws.Range("J23").Value = 123.3664479850047
$ws.Range("K23").Value = 876.3932961895841
$ws.Range("L23").Value = 539.9027985186409
$ws.Range("M23").Value = 1.827741282184856
$ws.Range("F24").Value = 12.68095976757514
$ws.Range("G24").Value = 3732
$ws.Range("H24").Value = 7.251116866025011
$ws.Range("I24").Value = 1.748828491096532
$ws.Range("J24").Value = 1.748828491096532
$ws.Range("K24").Value = 12.68095976757514
$ws.Range("L24").Value = 7.251116866025011
$ws.Range("M24").Value = 3.87656940094772
$ws.Range("J25").Value = 1.41203760553483
$ws.Range("K25").Value = 8.952846746409122
$ws.Range("L25").Value = 6.34037415952396
$ws.Range("M25").Value = 46.81788661726878
$ws.Range("N25").Value = 2134
$ws.Range("P25").Value = 0.0230776890906785

$ws.Range("Q25").Value = "(44.700207300046266, 48.93556593449126)"
$ws.Range("R25").Value = "(44.03465094320491, 49.601122291332615)"

Write-Host "Applied all changes"